$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("2.070225997297115e-07", "0.04240448674262143", "0.8054896365839992", "8.660232485948974", "1", "9.508126816298194")
    3 = @("3.230985683306322", "1.667794583268128", "0.8054896365839992", "0.496779210170732", "1", "6.201049113329182")
    4 = @("3.230985683306322", "1.667794583268128", "0.8054896365839992", "0.496779210170732", "1", "6.201049113329182")
    5 = @("1.459612070389937", "1.667794583268128", "0.8054896365839992", "0.496779210170732", "1", "4.429675500412797")
    6 = @("3.230985683306322", "1.667794583268128", "3.900430680208489", "0.496779210170732", "1", "9.295990156953671")
    7 = @("0.01514828764759746", "0.3127903958511391", "0.1575252929769615", "8.660232485948974", "0", "9.145696462424672")
    8 = @("3.230985683306322", "1.667794583268128", "0.1575252929769615", "0.496779210170732", "0", "5.553084769722144")
    9 = @("3.230985683306322", "1.667794583268128", "3.900430680208489", "0.496779210170732", "0", "9.295990156953671")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = [double]$vals[0]
    $ws.Cells.Item($row, 3).Value = [double]$vals[1]
    $ws.Cells.Item($row, 4).Value = [double]$vals[2]
    $ws.Cells.Item($row, 5).Value = [double]$vals[3]
    $ws.Cells.Item($row, 6).Value = [double]$vals[4]
    $ws.Cells.Item($row, 7).Value = [double]$vals[5]
}
